$wb = $excel.ActiveWorkbook

# Rename the 4th sheet (was "Sheet1") to "optionOnColumn"
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "optionOnColumn"

# Add a new worksheet right after "optionOnColumn" and name it "optionForTableAndColumn"
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "optionForTableAndColumn"

# Populate the new sheet with its table/column option test data.
# The order of assignment below matches the order the shared strings
# table should grow in (table-by-table: header-cell, then the "whole
# table" row, then any second column header for that table).
$ws5.Range("A1").Value = "####"
$ws5.Range("C1").Value = "optionForTableAndColumn001?opt1=val1#col1"
$ws5.Range("A3").Value = "optionForTableAndColumn001#~"

$ws5.Range("E1").Value = "optionForTableAndColumn002#col1?opt1=val1"
$ws5.Range("A4").Value = "optionForTableAndColumn002#~"
$ws5.Range("F1").Value = "optionForTableAndColumn002#col2?opt2=val2"

$ws5.Range("C3").Value = "aaa"
$ws5.Range("E4").Value = "bbb"
$ws5.Range("F4").Value = "ccc"

$ws5.Range("H1").Value = "optionForTableAndColumn003#col1"
$ws5.Range("A5").Value = "optionForTableAndColumn003?opt3=val3#~"
$ws5.Range("H5").Value = "ddd"

# Match column widths used in sibling sheets (values chosen so the
# engine's internal pixel rounding lands as close as possible to the
# widths recorded in the target workbook)
$ws5.Columns.Item(1).ColumnWidth = 39.571428571428573
$ws5.Columns.Item(2).ColumnWidth = 1.1428571428571428
$ws5.Columns.Item(3).ColumnWidth = 38.714285714285715
$ws5.Columns.Item(4).ColumnWidth = 1.1428571428571428
$ws5.Columns.Item(5).ColumnWidth = 41.857142857142857
$ws5.Columns.Item(6).ColumnWidth = 41.857142857142857
$ws5.Columns.Item(7).ColumnWidth = 1.1428571428571428
$ws5.Columns.Item(8).ColumnWidth = 31.857142857142854

# Make the new sheet active, with C1 selected/focused (matches "tabSelected" tab)
$ws5.Activate()
$ws5.Range("C1").Select()
